# Add bearings to BOM
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A45").Value = 3
$ws.Range("B45").Value = "Radial ball bearing, 4 mm x 9 mm x 4 mm (bore diameter x outer diameter x width)"
$ws.Range("D45").Value = "638-4-ZZ-SKF"
